$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update text/link/percentage cells (safe to assign directly) ---
$ws.Range('D2').Value = '67.681.34'
$ws.Range('E2').Value = '  -0.87%  '
$ws.Range('D3').Value = '3.768.29'
$ws.Range('E3').Value = '  -1.34%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('E5').Value = '  -0.70%  '
$ws.Range('E6').Value = '  +1.04%  '
$ws.Range('D7').Value = '3.767.48'
$ws.Range('E7').Value = '  -1.41%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('E10').Value = '  +1.36%  '
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('E12').Value = '  -0.54%  '
$ws.Range('E13').Value = '  +6.28%  '
$ws.Range('E14').Value = '  -0.72%  '
$ws.Range('D15').Value = '4.401.71'
$ws.Range('E15').Value = '  -1.30%  '
$ws.Range('D16').Value = '3.768.57'
$ws.Range('E16').Value = '  -1.28%  '
$ws.Range('E17').Value = '  +0.75%  '
$ws.Range('D18').Value = '67.629.30'
$ws.Range('E18').Value = '  -0.92%  '
$ws.Range('E19').Value = '  -2.48%  '
$ws.Range('E20').Value = '  +1.07%  '
$ws.Range('E21').Value = '  -4.89%  '
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('E23').Value = '  -1.42%  '
$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('E24').Value = '  -7.73%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('E25').Value = '  +1.23%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('E28').Value = '  +3.86%  '
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('E30').Value = '  -1.37%  '
$ws.Range('D31').Value = '3.913.50'
$ws.Range('E31').Value = '  -1.37%  '
$ws.Range('E32').Value = '  +0.11%  '
$ws.Range('E33').Value = '  -2.62%  '
$ws.Range('E34').Value = '  -2.74%  '
$ws.Range('E35').Value = '  -2.96%  '
$ws.Range('D36').Value = '3.732.66'
$ws.Range('E36').Value = '  -1.29%  '
$ws.Range('E37').Value = '  +5.30%  '
$ws.Range('E38').Value = '  +0.26%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('E39').Value = '  -1.61%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('E40').Value = '  -0.98%  '
$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('E41').Value = '  -1.12%  '
$ws.Range('E42').Value = '  -0.18%  '
$ws.Range('E43').Value = '  -0.39%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('E46').Value = '  -1.76%  '
$ws.Range('E47').Value = '  -2.13%  '
$ws.Range('E48').Value = '  -3.80%  '
$ws.Range('E49').Value = '  -8.05%  '
$ws.Range('E50').Value = '  -1.18%  '
$ws.Range('E51').Value = '  -0.82%  '

# --- Update numeric-looking price cells while preserving them as TEXT ---
# (set text format first so Excel does not auto-convert the string to a number)
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.16'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.19'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.526'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.51'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.455'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000278'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.76'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.67'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.57'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '469.49'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.722'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000148'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.04'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.23'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.19'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.40'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.68'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '30.66'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.25'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.16'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.88'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.138'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.87'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.998'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.313'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.77'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '45.96'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '401.05'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.000271'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '139.93'
$ws.Range('D50').Style = 'Normal'

Write-Host "Applied cryptos.xlsx update"
